$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation (2021-10-27, serial 44496) was recorded for
# this market/product. Insert it as the new first data row (row 2), pushing
# the existing rows down by one -- the oldest observation (formerly row 11)
# ends up at row 12.
$ws.Rows.Item(2).Insert()

# The inserted row picks up formatting from the row above (the bold header),
# so clear that before applying the real (unstyled) data formatting.
$ws.Range("A2:R2").ClearFormats()

$ws.Range("A2").Value = 10
$ws.Range("B2").Value = "Vega Modelo de Temuco"
$ws.Range("C2").Value = "La Araucanía"
$ws.Range("D2").Value = 44496
$ws.Range("D2").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E2").Value = 9
$ws.Range("F2").Value = 100112042
$ws.Range("G2").Value = "Locoto"
$ws.Range("H2").Value = "Sin especificar"
$ws.Range("I2").Value = "Primera"
$ws.Range("J2").Value = 40
$ws.Range("K2").Value = 2200
$ws.Range("L2").Value = 2200
$ws.Range("M2").Value = 2200
$ws.Range("N2").Value = "$/kilo"
$ws.Range("O2").Value = "Región de Arica y Parinacota"
$ws.Range("P2").Value = 2200
$ws.Range("Q2").Value = 1
$ws.Range("R2").Value = "Hortaliza"
